$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.949.17'
$ws.Range('E2').Value = '  +4.18%  '

$ws.Range('D3').Value = '2.778.55'
$ws.Range('E3').Value = '  +4.61%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.18'
$ws.Range('E5').Value = '  +0.33%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '161.19'
$ws.Range('E6').Value = '  +11.51%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  +3.67%  '

$ws.Range('E8').Value = '  -0.22%  '

$ws.Range('D9').Value = '2.803.81'
$ws.Range('E9').Value = '  +4.93%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.84'
$ws.Range('E10').Value = '  +3.45%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.114'

$ws.Range('E12').Value = '  +4.32%  '

$ws.Range('E13').Value = '  +1.27%  '

$ws.Range('D14').Value = '3.267.55'
$ws.Range('E14').Value = '  +4.47%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.81'
$ws.Range('E15').Value = '  +6.39%  '

$ws.Range('D16').Value = '63.889.08'
$ws.Range('E16').Value = '  +4.27%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000161'
$ws.Range('E17').Value = '  +9.32%  '

$ws.Range('D18').Value = '2.794.44'
$ws.Range('E18').Value = '  +4.98%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  +5.71%  '

$ws.Range('E20').Value = '  +5.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '367.86'
$ws.Range('E21').Value = '  +3.36%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.08'
$ws.Range('E22').Value = '  +2.48%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.549'
$ws.Range('E23').Value = '  +4.84%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('E25').Value = '  +4.70%  '

$ws.Range('E26').Value = '  +6.21%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.71'
$ws.Range('E27').Value = '  +2.54%  '

$ws.Range('D28').Value = '0.0₃0969'
$ws.Range('E28').Value = '  +17.40%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.35%  '

$ws.Range('E30').Value = '  +1.83%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.35'
$ws.Range('E31').Value = '  +5.57%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.26'
$ws.Range('E32').Value = '  +11.35%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '173.41'
$ws.Range('E33').Value = '  +2.34%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.94'
$ws.Range('E34').Value = '  +3.83%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.09'
$ws.Range('E35').Value = '  +8.73%  '

$ws.Range('E37').Value = '  +8.15%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.84'
$ws.Range('E38').Value = '  +7.42%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.04'
$ws.Range('E39').Value = '  +3.30%  '

$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.28'
$ws.Range('E40').Value = '  +3.14%  '

$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '341.90'
$ws.Range('E41').Value = '  -1.18%  '

$ws.Range('E42').Value = '  +15.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '39.73'
$ws.Range('E43').Value = '  +3.16%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '22.78'
$ws.Range('E44').Value = '  +9.41%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.84'
$ws.Range('E45').Value = '  +7.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0613'
$ws.Range('E46').Value = '  +5.90%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.655'
$ws.Range('E47').Value = '  +4.02%  '

$ws.Range('E48').Value = '  +3.15%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '138.60'
$ws.Range('E49').Value = '  +2.11%  '

$ws.Range('D51').Value = '2.185.42'
$ws.Range('E51').Value = '  +4.01%  '
